$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column cells that look numeric stay as text
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D10", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D24", "D25", "D26", "D27", "D29", "D32", "D33", "D34", "D36", "D37", "D40", "D41", "D43", "D44", "D46", "D48", "D49", "D50", "D51")
foreach ($pc in $priceCells) { $ws.Range($pc).NumberFormat = "@" }

$ws.Range('D2').Value = '37.709.08'
$ws.Range('E2').Value = '  -1.41%  '
$ws.Range('D3').Value = '2.028.61'
$ws.Range('E3').Value = '  -1.77%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '227.53'
$ws.Range('D6').Value = '0.607'
$ws.Range('E6').Value = '  -1.17%  '
$ws.Range('D7').Value = '59.79'
$ws.Range('E7').Value = '  -2.32%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  -2.95%  '
$ws.Range('D10').Value = '0.0825'
$ws.Range('E10').Value = '  +2.20%  '
$ws.Range('E11').Value = '  -0.26%  '
$ws.Range('D12').Value = '14.62'
$ws.Range('E12').Value = '  -1.66%  '
$ws.Range('D13').Value = '2.329.54'
$ws.Range('E13').Value = '  -1.72%  '
$ws.Range('D14').Value = '20.97'
$ws.Range('E14').Value = '  -1.38%  '
$ws.Range('D15').Value = '0.765'
$ws.Range('E15').Value = '  +0.13%  '
$ws.Range('D16').Value = '5.18'
$ws.Range('E16').Value = '  -2.67%  '
$ws.Range('D17').Value = '2.028.76'
$ws.Range('E17').Value = '  -2.02%  '
$ws.Range('D18').Value = '37.679.23'
$ws.Range('E18').Value = '  -1.38%  '
$ws.Range('D19').Value = '69.48'
$ws.Range('E19').Value = '  -0.98%  '
$ws.Range('D20').Value = '5.88'
$ws.Range('E20').Value = '  -6.25%  '
$ws.Range('D21').Value = '0.0₃0823'
$ws.Range('E21').Value = '  -1.63%  '
$ws.Range('D22').Value = '223.45'
$ws.Range('E22').Value = '  -1.18%  '
$ws.Range('D24').Value = '2.34'
$ws.Range('E24').Value = '  -3.00%  '
$ws.Range('D25').Value = '2.28'
$ws.Range('E25').Value = '  +2.09%  '
$ws.Range('D26').Value = '167.49'
$ws.Range('E26').Value = '  +0.73%  '
$ws.Range('D27').Value = '9.34'
$ws.Range('E27').Value = '  +0.96%  '
$ws.Range('E28').Value = '  -3.18%  '
$ws.Range('D29').Value = '18.74'
$ws.Range('E29').Value = '  -1.24%  '
$ws.Range('E30').Value = '  -3.81%  '
$ws.Range('E31').Value = '  +0.63%  '
$ws.Range('D32').Value = '2.22'
$ws.Range('E32').Value = '  +8.55%  '
$ws.Range('D33').Value = '4.37'
$ws.Range('E33').Value = '  -4.30%  '
$ws.Range('D34').Value = '0.0605'
$ws.Range('E34').Value = '  -0.07%  '
$ws.Range('E35').Value = '  -2.98%  '
$ws.Range('D36').Value = '6.43'
$ws.Range('E36').Value = '  +2.56%  '
$ws.Range('D37').Value = '2.31'
$ws.Range('E37').Value = '  -0.62%  '
$ws.Range('E38').Value = '  +2.92%  '
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('D40').Value = '18.10'
$ws.Range('E40').Value = '  +6.13%  '
$ws.Range('D41').Value = '1.535.17'
$ws.Range('E41').Value = '  +0.87%  '
$ws.Range('E42').Value = '  -0.88%  '
$ws.Range('D43').Value = '95.73'
$ws.Range('E43').Value = '  -2.35%  '
$ws.Range('D44').Value = '2.79'
$ws.Range('E44').Value = '  -2.61%  '
$ws.Range('E45').Value = '  -2.25%  '
$ws.Range('D46').Value = '4.06'
$ws.Range('E46').Value = '  +0.71%  '
$ws.Range('E47').Value = '  -2.80%  '
$ws.Range('B48').Value = 'MXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D48').Value = '2.96'
$ws.Range('E48').Value = '  -0.15%  '
$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').Value = '1.00'
$ws.Range('E49').Value = '  -1.91%  '
$ws.Range('D50').Value = '7.08'
$ws.Range('E50').Value = '  -0.51%  '
$ws.Range('D51').Value = '2.219.78'
$ws.Range('E51').Value = '  -1.69%  '
